$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "26.792.97"
$c.Style = "Normal"
$ws.Range("E2").Value = "  +1.14%  "
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "1.648.99"
$c.Style = "Normal"
$ws.Range("E3").Value = "  +1.38%  "
$ws.Range("E4").Value = "  +0.71%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "216.55"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +1.79%  "
$ws.Range("E6").Value = "  +1.05%  "
$ws.Range("E7").Value = "  +0.61%  "
$ws.Range("E8").Value = "  +1.58%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.0628"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +2.23%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.0844"
$c.Style = "Normal"
$ws.Range("E11").Value = "  +0.19%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "1.879.20"
$c.Style = "Normal"
$ws.Range("E12").Value = "  +1.46%  "
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "1.650.20"
$c.Style = "Normal"
$ws.Range("E13").Value = "  +1.03%  "
$ws.Range("E14").Value = "  +1.55%  "
$ws.Range("E15").Value = "  +2.04%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "65.52"
$c.Style = "Normal"
$ws.Range("E16").Value = "  +0.90%  "
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "26.810.14"
$c.Style = "Normal"
$ws.Range("E17").Value = "  +1.16%  "
$ws.Range("E18").Value = "  +0.78%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "218.59"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +1.97%  "
$ws.Range("E20").Value = "  +0.57%  "
$ws.Range("E21").Value = "  +1.90%  "
$ws.Range("E22").Value = "  +16.78%  "
$ws.Range("E23").Value = "  +0.39%  "
$ws.Range("E24").Value = "  +2.62%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "146.73"
$c.Style = "Normal"
$ws.Range("E25").Value = "  -1.26%  "
$ws.Range("E26").Value = "  +0.67%  "
$ws.Range("E27").Value = "  +0.31%  "
$ws.Range("E28").Value = "  +3.94%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "15.73"
$c.Style = "Normal"
$ws.Range("E29").Value = "  +1.30%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "0.0518"
$c.Style = "Normal"
$ws.Range("E30").Value = "  +1.98%  "
$ws.Range("E31").Value = "  +2.27%  "
$ws.Range("E32").Value = "  +0.67%  "
$ws.Range("E33").Value = "  +2.11%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "1.285.87"
$c.Style = "Normal"
$ws.Range("E34").Value = "  +3.90%  "
$ws.Range("E35").Value = "  +3.24%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "2.45"
$c.Style = "Normal"
$ws.Range("E36").Value = "  +3.03%  "
$ws.Range("E37").Value = "  +3.19%  "
$ws.Range("E38").Value = "  +6.14%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "0.829"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +0.65%  "
$ws.Range("E41").Value = "  +2.23%  "
$ws.Range("E42").Value = "  -0.70%  "
$ws.Range("E43").Value = "  +2.60%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "1.789.80"
$c.Style = "Normal"
$ws.Range("E44").Value = "  +1.61%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "92.07"
$c.Style = "Normal"
$ws.Range("E45").Value = "  -0.91%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "59.74"
$c.Style = "Normal"
$ws.Range("E47").Value = "  +1.82%  "
$ws.Range("E48").Value = "  +1.47%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "7.76"
$c.Style = "Normal"
$ws.Range("E49").Value = "  +3.91%  "
$ws.Range("E50").Value = "  +2.20%  "
$ws.Range("E51").Value = "  +0.51%  "
